$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"; $ws.Range("D2").Value = '45.960.03'; $ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"; $ws.Range("E2").Value = '  -1.59%  '; $ws.Range("E2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"; $ws.Range("D3").Value = '2.378.13'; $ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"; $ws.Range("E3").Value = '  +2.92%  '; $ws.Range("E3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"; $ws.Range("D5").Value = '300.56'; $ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"; $ws.Range("E5").Value = '  -0.46%  '; $ws.Range("E5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"; $ws.Range("D6").Value = '98.97'; $ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"; $ws.Range("E6").Value = '  -3.69%  '; $ws.Range("E6").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"; $ws.Range("E7").Value = '  -0.79%  '; $ws.Range("E7").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"; $ws.Range("E8").Value = '  +0.03%  '; $ws.Range("E8").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"; $ws.Range("E9").Value = '  -3.83%  '; $ws.Range("E9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"; $ws.Range("D10").Value = '34.55'; $ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"; $ws.Range("E10").Value = '  -6.54%  '; $ws.Range("E10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"; $ws.Range("D11").Value = '0.0791'; $ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"; $ws.Range("E11").Value = '  -1.92%  '; $ws.Range("E11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"; $ws.Range("D12").Value = '7.16'; $ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"; $ws.Range("E12").Value = '  -3.29%  '; $ws.Range("E12").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"; $ws.Range("E13").Value = '  -0.48%  '; $ws.Range("E13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"; $ws.Range("D14").Value = '2.737.52'; $ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"; $ws.Range("E14").Value = '  +2.74%  '; $ws.Range("E14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"; $ws.Range("D15").Value = '2.374.57'; $ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"; $ws.Range("E15").Value = '  +2.76%  '; $ws.Range("E15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"; $ws.Range("D16").Value = '0.817'; $ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"; $ws.Range("E16").Value = '  -0.74%  '; $ws.Range("E16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"; $ws.Range("D17").Value = '13.67'; $ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"; $ws.Range("E17").Value = '  -2.85%  '; $ws.Range("E17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"; $ws.Range("D18").Value = '45.870.26'; $ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"; $ws.Range("E18").Value = '  -1.75%  '; $ws.Range("E18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"; $ws.Range("D19").Value = '12.75'; $ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"; $ws.Range("E19").Value = '  -3.72%  '; $ws.Range("E19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"; $ws.Range("D20").Value = '0.0₃0957'; $ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"; $ws.Range("E20").Value = '  +0.90%  '; $ws.Range("E20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"; $ws.Range("D21").Value = '6.03'; $ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"; $ws.Range("E21").Value = '  -2.01%  '; $ws.Range("E21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"; $ws.Range("D22").Value = '67.29'; $ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"; $ws.Range("E22").Value = '  +0.44%  '; $ws.Range("E22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"; $ws.Range("D23").Value = '244.42'; $ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"; $ws.Range("E23").Value = '  -1.64%  '; $ws.Range("E23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"; $ws.Range("D24").Value = '2.80'; $ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"; $ws.Range("E24").Value = '  -5.06%  '; $ws.Range("E24").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"; $ws.Range("E25").Value = '  +0.01%  '; $ws.Range("E25").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"; $ws.Range("E26").Value = '  -2.24%  '; $ws.Range("E26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"; $ws.Range("D27").Value = '39.49'; $ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"; $ws.Range("E27").Value = '  -9.63%  '; $ws.Range("E27").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"; $ws.Range("E28").Value = '  -3.40%  '; $ws.Range("E28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"; $ws.Range("D29").Value = '9.77'; $ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"; $ws.Range("E29").Value = '  -2.22%  '; $ws.Range("E29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"; $ws.Range("D30").Value = '3.81'; $ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"; $ws.Range("E30").Value = '  +21.85%  '; $ws.Range("E30").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"; $ws.Range("E31").Value = '  +4.36%  '; $ws.Range("E31").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"; $ws.Range("E32").Value = '  +6.90%  '; $ws.Range("E32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"; $ws.Range("D33").Value = '5.54'; $ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"; $ws.Range("E33").Value = '  -4.55%  '; $ws.Range("E33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"; $ws.Range("D34").Value = '146.88'; $ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"; $ws.Range("E34").Value = '  +0.52%  '; $ws.Range("E34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"; $ws.Range("D35").Value = '0.0774'; $ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"; $ws.Range("E35").Value = '  -3.99%  '; $ws.Range("E35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"; $ws.Range("D36").Value = '0.112'; $ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"; $ws.Range("E36").Value = '  +0.15%  '; $ws.Range("E36").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"; $ws.Range("E37").Value = '  +7.43%  '; $ws.Range("E37").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"; $ws.Range("E38").Value = '  -3.32%  '; $ws.Range("E38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"; $ws.Range("D39").Value = '14.89'; $ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"; $ws.Range("E39").Value = '  -5.20%  '; $ws.Range("E39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"; $ws.Range("D40").Value = '3.90'; $ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"; $ws.Range("E40").Value = '  -5.55%  '; $ws.Range("E40").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"; $ws.Range("E41").Value = '  -2.12%  '; $ws.Range("E41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"; $ws.Range("D42").Value = '3.22'; $ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"; $ws.Range("E42").Value = '  -7.22%  '; $ws.Range("E42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"; $ws.Range("D43").Value = '1.934.56'; $ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"; $ws.Range("E43").Value = '  +4.31%  '; $ws.Range("E43").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"; $ws.Range("D45").Value = '91.88'; $ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"; $ws.Range("E45").Value = '  +2.91%  '; $ws.Range("E45").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"; $ws.Range("E46").Value = '  -9.85%  '; $ws.Range("E46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"; $ws.Range("D47").Value = '8.45'; $ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"; $ws.Range("E47").Value = '  +5.43%  '; $ws.Range("E47").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"; $ws.Range("E48").Value = '  -5.38%  '; $ws.Range("E48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"; $ws.Range("D49").Value = '2.608.78'; $ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"; $ws.Range("E49").Value = '  +2.64%  '; $ws.Range("E49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"; $ws.Range("D50").Value = '97.86'; $ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"; $ws.Range("E50").Value = '  +0.29%  '; $ws.Range("E50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"; $ws.Range("D51").Value = '68.62'; $ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"; $ws.Range("E51").Value = '  -8.25%  '; $ws.Range("E51").Style = "Normal"
